$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column cells are forced back to Text (matching the workbook's existing
# inlineStr convention) since some prices parse as valid numbers and would
# otherwise be auto-converted to the Number type by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.743.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.96%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.987.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.83%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '497.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.70%  '

$ws.Range("E8").Value = '  -3.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.106'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.73%  '

$ws.Range("E11").Value = '  -5.06%  '

$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.494.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.680.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000149'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.983.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.51%  '

$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.96%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  -7.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.79%  '

$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("E26").Value = '  -2.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0893'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("E31").Value = '  -5.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.46%  '

$ws.Range("E33").Value = '  -7.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.21%  '

$ws.Range("E36").Value = '  -5.34%  '

$ws.Range("E37").Value = '  -8.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0673'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.018.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.36'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("E43").Value = '  -6.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.234.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.34%  '

$ws.Range("E47").Value = '  -8.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.77%  '

$ws.Range("E49").Value = '  +1.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.57%  '

$ws.Range("E51").Value = '  -5.67%  '
